$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# User Story cell: TC-134 -> TC-63812
$ws.Range("B4").Value = "TC-63812"

# Row Index column (E) updated from 13 to 6 for all data rows
$ws.Range("E8").Value = 6
$ws.Range("E9").Value = 6
$ws.Range("E10").Value = 6
$ws.Range("E11").Value = 6
$ws.Range("E12").Value = 6
$ws.Range("E13").Value = 6
$ws.Range("E14").Value = 6

# Additonal Base Row Index column (G) updated from 13 to 6 (only rows with a
# numeric Additional Base Row Index; rows 11-14 hold "NA" and stay untouched)
$ws.Range("G8").Value = 6
$ws.Range("G9").Value = 6
$ws.Range("G10").Value = 6

# Update the active selection shown when the workbook is reopened
$ws.Range("B5").Select()
